$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1949.5
$ws.Range("I80").Value = 2549.5
$ws.Range("K80").Value = 7648.5
$ws.Range("M80").Value = -6650.5

$ws.Range("H83").Value = 1949.5
$ws.Range("I83").Value = 2549.5
$ws.Range("K83").Value = 22945.5
$ws.Range("M83").Value = -17953.5

$ws.Range("H98").Value = 6316
$ws.Range("I98").Value = 6316
$ws.Range("K98").Value = 6316
$ws.Range("M98").Value = -4818

$ws.Range("H106").Value = 1237.25
$ws.Range("I106").Value = 983
$ws.Range("J106").Value = 2000
$ws.Range("K106").Value = 983
$ws.Range("L106").Value = 2000
$ws.Range("M106").Value = -352
$ws.Range("N106").Value = -3262

$ws.Range("H122").Value = 6316
$ws.Range("I122").Value = 6316
$ws.Range("K122").Value = 18948
$ws.Range("M122").Value = -16498

$ws.Range("H137").Value = 2870.5715
$ws.Range("I137").Value = 2682.1667
$ws.Range("K137").Value = 8046.500100000001
$ws.Range("M137").Value = -5496.500100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 310500.38
$ws.Range("I74").Value = 310500.38
$ws.Range("K74").Value = 310500.38
$ws.Range("M74").Value = -309626.38

$ws.Range("H77").Value = 310500.38
$ws.Range("I77").Value = 310500.38
$ws.Range("K77").Value = 1552501.9
$ws.Range("M77").Value = -1548133.9

$ws.Range("H132").Value = 3608.2
$ws.Range("J132").Value = 3192.5
$ws.Range("L132").Value = 9577.5
$ws.Range("N132").Value = -14637.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1463.8572
$ws.Range("I20").Value = 1294.8
$ws.Range("J20").Value = 1886.5
$ws.Range("K20").Value = 1294.8
$ws.Range("L20").Value = 1886.5
$ws.Range("M20").Value = -1047.8
$ws.Range("N20").Value = -2380.5

$ws.Range("H99").Value = 2775
$ws.Range("I99").Value = 2775
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2775
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1277
$ws.Range("N99").Value = $null

$ws.Range("H105").Value = 2449.5
$ws.Range("I105").Value = 2449.5
$ws.Range("K105").Value = 2449.5
$ws.Range("M105").Value = -702.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 4154.143
$ws.Range("I35").Value = 3804.8
$ws.Range("J35").Value = 5027.5
$ws.Range("K35").Value = 3804.8
$ws.Range("L35").Value = 5027.5
$ws.Range("M35").Value = -3510.8
$ws.Range("N35").Value = -5615.5

$ws.Range("H58").Value = 398.05884
$ws.Range("I58").Value = 398.05884
$ws.Range("K58").Value = 398.05884
$ws.Range("M58").Value = -195.05884

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").Value = $null

$ws.Range("H132").Value = 6182.375
$ws.Range("I132").Value = 5666.3335
$ws.Range("K132").Value = 16999.0005
$ws.Range("M132").Value = -14469.0005

$ws.Range("H136").Value = 398.05884
$ws.Range("I136").Value = 398.05884
$ws.Range("K136").Value = 1194.17652
$ws.Range("M136").Value = 1355.82348

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 68.64286
$ws.Range("I2").Value = 38.77778
$ws.Range("K2").Value = 232.66668
$ws.Range("M2").Value = -119.66668

$ws.Range("H4").Value = 77231.80499999999
$ws.Range("I4").Value = 41121.746
$ws.Range("J4").Value = 501525
$ws.Range("K4").Value = 123365.238
$ws.Range("L4").Value = 1504575
$ws.Range("M4").Value = -123253.238
$ws.Range("N4").Value = -1504799

$ws.Range("H6").Value = 29
$ws.Range("J6").Value = 1.5
$ws.Range("L6").Value = 4.5
$ws.Range("N6").Value = -230.5

$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = $null

$ws.Range("H38").Value = 3123.75
$ws.Range("I38").Value = 2497.5
$ws.Range("K38").Value = 7492.5
$ws.Range("M38").Value = -7145.5

$ws.Range("H108").Value = 27
$ws.Range("I108").Value = 27
$ws.Range("K108").Value = 81
$ws.Range("M108").Value = 2799

$ws.Range("H131").Value = 762.8889
$ws.Range("I131").Value = 733.25
$ws.Range("K131").Value = 2199.75
$ws.Range("M131").Value = 2840.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1083.3334
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 1325
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 1325
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -1915

$ws.Range("H27").Value = 1083.3334
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 1325
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 1325
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -1539

$ws.Range("H40").Value = 320557.2
$ws.Range("I40").Value = 9746.223
$ws.Range("J40").Value = 720171.3
$ws.Range("K40").Value = 9746.223
$ws.Range("L40").Value = 720171.3
$ws.Range("M40").Value = -9610.223
$ws.Range("N40").Value = -720443.3

$ws.Range("H46").Value = 592.3333
$ws.Range("I46").Value = 592.3333
$ws.Range("K46").Value = 592.3333
$ws.Range("M46").Value = -404.3333

$ws.Range("H136").Value = 867999.8
$ws.Range("I136").Value = 1698333
$ws.Range("J136").Value = 37666.668
$ws.Range("K136").Value = 5094999
$ws.Range("L136").Value = 113000.004
$ws.Range("M136").Value = -5092449
$ws.Range("N136").Value = -118100.004

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 530.3333
$ws.Range("I107").Value = 239.7
$ws.Range("J107").Value = 1983.5
$ws.Range("K107").Value = 719.0999999999999
$ws.Range("L107").Value = 5950.5
$ws.Range("M107").Value = 1200.9
$ws.Range("N107").Value = -9790.5
